{"js": "// 1) Locate the paragraphs whose full text matches the PERSON_75/76/101 list items.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nlet p75 = null;\nlet p76 = null;\nlet p101 = null;\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t === \"[[PERSON_75]] \u2013 [[PERSON_75]], [[PERSON_75]]\") {\n    p75 = items[i];\n  } else if (t === \"[[PERSON_76]] \u2013 [[PERSON_76]], [[PERSON_76]]\") {\n    p76 = items[i];\n  } else if (t === \"[[PERSON_101]] \u2013 [[PERSON_101]], [[PERSON_101]]\") {\n    p101 = items[i];\n  }\n}\n\nif (!p75 || !p76 || !p101) {\n  throw new Error(\"Could not locate all target paragraphs: \" +\n    \"p75=\" + !!p75 + \" p76=\" + !!p76 + \" p101=\" + !!p101);\n}\n\n// 2) Rewrite PERSON_75's paragraph text so the middle token becomes PERSON_76.\np75.insertText(\"[[PERSON_75]] \u2013 [[PERSON_76]], [[PERSON_75]]\", \"Replace\");\n\n// 3) Insert the replacement list item (PERSON_102) right after PERSON_101's paragraph,\n//    inheriting the same numbered-list formatting.\np101.insertParagraph(\"[[PERSON_102]] \u2013 [[PERSON_102]], [[PERSON_102]]\", \"After\");\n\n// 4) Remove the old PERSON_76 paragraph entirely.\np76.delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$dash = [char]0x2013\n\n$target75 = \"[[PERSON_75]] \" + $dash + \" [[PERSON_75]], [[PERSON_75]]\"\n$target76 = \"[[PERSON_76]] \" + $dash + \" [[PERSON_76]], [[PERSON_76]]\"\n$target101 = \"[[PERSON_101]] \" + $dash + \" [[PERSON_101]], [[PERSON_101]]\"\n\n$p75 = $null\n$p76 = $null\n$p101 = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]0x0D, [char]0x07)\n    if ($t -eq $target75) {\n        $p75 = $p\n    } elseif ($t -eq $target76) {\n        $p76 = $p\n    } elseif ($t -eq $target101) {\n        $p101 = $p\n    }\n}\n\nif ($p75 -eq $null -or $p76 -eq $null -or $p101 -eq $null) {\n    throw \"Could not locate all target paragraphs\"\n}\n\n# 1) Rewrite PERSON_75's paragraph text so the middle token becomes PERSON_76.\n$p75.Range.Text = \"[[PERSON_75]] \" + $dash + \" [[PERSON_76]], [[PERSON_75]]\"\n\n# 2) Insert a new list paragraph right after PERSON_101's paragraph, inheriting\n#    its numbered-list formatting, and set its text to the PERSON_102 variant.\n$p101.Range.InsertParagraphAfter()\n$newPara = $p101.Next()\n$newPara.Range.Text = \"[[PERSON_102]] \" + $dash + \" [[PERSON_102]], [[PERSON_102]]\"\n\n# 3) Remove the old PERSON_76 paragraph entirely.\n$p76.Range.Delete()\n"}
